$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the form: drop the unused trailing columns (K:V) and
# rewrite the header row with the new, shorter set of fields. ---
$ws.Columns("K:V").Delete()

$ws.Range("A1").Value = "nim"
$ws.Range("B1").Value = "nama"
$ws.Range("C1").Value = "nama_sma"
$ws.Range("D1").Value = "penghasilan_orang_tua"
$ws.Range("E1").Value = "prodi"
$ws.Range("F1").Value = "jalur"
$ws.Range("G1").Value = "ip_semester_1"
$ws.Range("H1").Value = "ip_semester_2"
$ws.Range("I1").Value = "ip_semester_3"
$ws.Range("J1").Value = "ip_semester_4"

# --- Column widths (best-fit to the new header text) ---
$ws.Columns(1).ColumnWidth = 3.6666666666666665
$ws.Columns(2).ColumnWidth = 5.333333333333333
$ws.Columns(3).ColumnWidth = 10
$ws.Columns(4).ColumnWidth = 21.333333333333332
$ws.Columns(5).ColumnWidth = 4.833333333333333
$ws.Columns(6).ColumnWidth = 4.333333333333333
$ws.Columns(7).ColumnWidth = 12.666666666666666
$ws.Columns(8).ColumnWidth = 12.666666666666666
$ws.Columns(9).ColumnWidth = 12.666666666666666
$ws.Columns(10).ColumnWidth = 12.666666666666666

# --- Two blank placeholder rows underneath the header for data entry,
# styled with the green monospace "input" font used elsewhere in the
# workbook. ---
$rng1 = $ws.Range("C2")
$rng1.Font.Name = "Consolas"
$rng1.Font.Family = 3
$rng1.Font.Color = 7979928
$rng1.VerticalAlignment = -4108

$rng1.Copy()
$ws.Range("C2:C3").PasteSpecial(-4122)
$ws.Range("E2:F3").PasteSpecial(-4122)

# --- Page setup / selection to match the refreshed UI state ---
$ws.PageSetup.Orientation = 1
$ws.Range("D6").Select() | Out-Null
